# Update header row labels so that Power BI can automatically treat the
# first row of each table as a header (prefixing the year/interval values
# with a descriptive word).

$wb = $excel.ActiveWorkbook

# Sheets whose header row is "2015 / 2030 / 2040 / 2050" -> "Ano 2015 / Ano 2030 / Ano 2040 / Ano 2050"
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# "Potencia Incremental - SIN(MW)" header row uses year intervals instead of single years
$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo 2015"
$wsIntervalo.Range("C1").Value = "Intervalo 2015-2030"
$wsIntervalo.Range("D1").Value = "Intervalo 2031-2040"
$wsIntervalo.Range("E1").Value = "Intervalo 2041-2050"

# "Custo Total (bilhões de R$)" only has a single year column in the header row
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
